$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.250.90"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").Value = "1.912.57"
$ws.Range("E3").Value = "  +2.22%  "

$ws.Range("D4").Value = "'0.9977"
$ws.Range("E4").Value = "  -0.33%  "

$ws.Range("D5").Value = "'313.32"
$ws.Range("E5").Value = "  +0.50%  "

$ws.Range("D6").Value = "'0.9981"
$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("D7").Value = "'0.5077"
$ws.Range("E7").Value = "  +0.39%  "

$ws.Range("D8").Value = "'0.3945"
$ws.Range("E8").Value = "  +0.47%  "

$ws.Range("D9").Value = "'0.09360"
$ws.Range("E9").Value = "  -3.35%  "

$ws.Range("D10").Value = "'1.143"
$ws.Range("E10").Value = "  +0.29%  "

$ws.Range("D11").Value = "'41.90"
$ws.Range("E11").Value = "  +2.46%  "

$ws.Range("D12").Value = "'6.395"
$ws.Range("E12").Value = "  -1.85%  "

$ws.Range("D13").Value = "'20.93"
$ws.Range("E13").Value = "  +0.03%  "

$ws.Range("D14").Value = "1.899.98"
$ws.Range("E14").Value = "  +1.09%  "

$ws.Range("D15").Value = "'7.323"
$ws.Range("E15").Value = "  -1.13%  "

$ws.Range("D16").Value = "'0.9976"
$ws.Range("E16").Value = "  -0.35%  "

$ws.Range("D17").Value = "'0.00001127"
$ws.Range("E17").Value = "  -0.17%  "

$ws.Range("D18").Value = "'92.69"
$ws.Range("E18").Value = "  -0.14%  "

$ws.Range("D19").Value = "'0.06580"
$ws.Range("E19").Value = "  -0.14%  "

$ws.Range("D20").Value = "'17.94"
$ws.Range("E20").Value = "  +2.16%  "

$ws.Range("D21").Value = "'0.9981"
$ws.Range("E21").Value = "  -0.24%  "

$ws.Range("D22").Value = "'6.227"
$ws.Range("E22").Value = "  +0.97%  "

$ws.Range("D23").Value = "28.285.14"
$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").Value = "'11.41"
$ws.Range("E24").Value = "  +0.44%  "

$ws.Range("D25").Value = "'2.309"
$ws.Range("E25").Value = "  +1.46%  "

$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.630"
$ws.Range("E26").Value = "  +3.78%  "

$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.119.47"
$ws.Range("E27").Value = "  +2.16%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'21.06"
$ws.Range("E28").Value = "  -0.91%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "'157.42"
$ws.Range("E29").Value = "  -0.52%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'127.36"
$ws.Range("E30").Value = "  -0.10%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.097"
$ws.Range("E31").Value = "  +2.72%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.1071"
$ws.Range("E32").Value = "  +0.87%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.638"
$ws.Range("E33").Value = "  +0.10%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'3.609"
$ws.Range("E34").Value = "  -0.51%  "

$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").Value = "'9.683"
$ws.Range("E35").Value = "  +1.15%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.06702"
$ws.Range("E36").Value = "  -0.13%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02418"
$ws.Range("E37").Value = "  +1.07%  "

$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'1.256"
$ws.Range("E38").Value = "  +1.40%  "

$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.2188"
$ws.Range("E39").Value = "  +0.38%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.252"
$ws.Range("E40").Value = "  +6.38%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.6410"
$ws.Range("E41").Value = "  +0.72%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'11.54"
$ws.Range("E42").Value = "  +0.36%  "

$ws.Range("B43").Value = "InternetComputer(DFINITY)"
$ws.Range("C43").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D43").Value = "'4.998"
$ws.Range("E43").Value = "  +0.66%  "

$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "'0.9985"
$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'13.37"
$ws.Range("E45").Value = "  -1.05%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.6011"
$ws.Range("E46").Value = "  +0.06%  "

$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.706"
$ws.Range("E47").Value = "  +1.32%  "

$ws.Range("B48").Value = "WEMIXTOKEN"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "'1.276"
$ws.Range("E48").Value = "  +1.48%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'2.022"
$ws.Range("E49").Value = "  +1.56%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'123.19"
$ws.Range("E50").Value = "  -0.69%  "

$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "'1.182"
$ws.Range("E51").Value = "  -1.13%  "
